$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Establish new shared strings in the same relative order the author created them in ---
# (PNG and "Image Decode (ms)" already exist in the workbook's string table from before)
$ws.Range("B17").Value = "LCP"
$ws.Range("B2").Value = "JPG"
$ws.Range("A14").Value = "Среднее"
$ws.Range("A15").Value = "Стандартное отклонение"
$ws.Range("D18").Value = "AVIF"
$ws.Range("E18").Value = "WebP"
$ws.Range("A26").Value = "Стандратное отклонение"

# --- Section 1: PNG Image Decode data, rewritten with new values/layout (rows 1-15) ---

# Headers
$ws.Range("B1").Value = "Image Decode (ms)"

# Data values B4:B13
$ws.Range("B4").Value = 12.51
$ws.Range("B5").Value = 15.76
$ws.Range("B6").Value = 13.4
$ws.Range("B7").Formula = "=7.44+5.73"
$ws.Range("B8").Value = 12.97
$ws.Range("B9").Value = 12.68
$ws.Range("B10").Value = 11.9
$ws.Range("B11").Value = 14.16
$ws.Range("B12").Value = 13.84
$ws.Range("B13").Value = 14.04

# Summary rows
$ws.Range("B14").Formula = "=AVERAGE(B4:B13)"
$ws.Range("B15").Formula = "=STDEV(B4:B13)"

# --- Section 2: Image comparison table (rows 17-26) ---

$ws.Range("B18").Value = "JPG"
$ws.Range("C18").Value = "PNG"

$ws.Range("B19").Value = 81.36
$ws.Range("C19").Value = 135.9
$ws.Range("D19").Value = 133.80000000000001
$ws.Range("E19").Value = 130.94

$ws.Range("B20").Value = 65.63
$ws.Range("C20").Value = 151.69999999999999
$ws.Range("D20").Value = 137
$ws.Range("E20").Value = 122.2

$ws.Range("B21").Value = 68.63
$ws.Range("C21").Value = 143.4
$ws.Range("D21").Value = 152.80000000000001
$ws.Range("E21").Value = 136.01

$ws.Range("B22").Value = 66.349999999999994
$ws.Range("C22").Value = 149.9
$ws.Range("D22").Value = 135.80000000000001
$ws.Range("E22").Value = 125.7

$ws.Range("B23").Value = 87.21
$ws.Range("C23").Value = 133.94999999999999
$ws.Range("D23").Value = 142.9
$ws.Range("E23").Value = 116.3

$ws.Range("B24").Value = 71.87
$ws.Range("C24").Value = 133
$ws.Range("D24").Value = 135
$ws.Range("E24").Value = 121.2

$ws.Range("A25").Value = "Среднее"
$ws.Range("B25").Formula = "=AVERAGE(B19:B24)"
$ws.Range("C25:E25").Formula = "=AVERAGE(C19:C24)"

$ws.Range("B26").Formula = "=STDEV(B19:B25)"
$ws.Range("C26:E26").Formula = "=STDEV(C19:C25)"

# --- Number formatting (style index 1 with numFmtId 2 "0.00") applied to all numeric-looking cells ---
$ws.Range("B4:B15").NumberFormat = "0.00"
$ws.Range("B19:E26").NumberFormat = "0.00"

# --- Column width for column B ---
$ws.Columns("B").ColumnWidth = 11.14

# --- Sheet view: selection & scroll position ---
$ws.Range("S45").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
